# NaukriSearch.xlsx edit
#
# The NaukriSearch tab tracked a single "current" requisition in row 2.
# This edit archives the existing "Java developer" requisition into
# Sheet1 (as a new row) and replaces row 2 on NaukriSearch with a fresh
# "Developer" requisition.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NaukriSearch")
$ws2 = $wb.Worksheets.Item("Sheet1")

# 1. Archive the current NaukriSearch requisition (row 2) into Sheet1 as a
#    new row 3, carrying over values + formatting before row 2 is overwritten.
$ws1.Range("A2:J2").Copy($ws2.Range("A3:J3"))
$ws2.Range("J3").ClearContents()
$ws2.Range("C3").Style = $ws2.Range("A3").Style
$ws2.Range("I3").Style = $ws2.Range("A3").Style

# 2. Replace the NaukriSearch requisition in row 2 with the new posting.
$ws1.Range("B2").Value = "Developer"
$ws1.Range("C2").Value = "Asp.net, C#"
$ws1.Range("D2").Value = "4-6"
$ws1.Range("E2").Value = "Bangalore"
$ws1.Range("H2").Value = 25

# 3. The SkillKeywords column needs to be a bit wider to fit the new text.
$ws1.Columns.Item(3).AutoFit()

# 4. Leave the UI selection the way the author left it: Sheet1 highlights
#    the newly archived row, and NaukriSearch (still the active tab) is
#    parked on H3. Do Sheet1 first so NaukriSearch ends up active/last.
$ws2.Activate()
$ws2.Rows.Item(3).Select()

$ws1.Activate()
$ws1.Range("H3").Select()
